# Update gh-pages to output generated at 456a3b4
# Updates the "想去人数" (F column) counts on the 展览 and 全部类型 sheets.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F2"  = 2912
    "F5"  = 6707
    "F6"  = 1643
    "F9"  = 53
    "F10" = 109
    "F11" = 22
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}
